# Swap the presentation's theme palette: the deck currently carries the
# "Integral" color scheme (green palette) on its live/active theme; the
# edit replaces it with the stock "Office Theme" palette (the colors that
# used to live in the sibling theme part), matching the authored change
# where the two theme parts' contents were exchanged.
#
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink -- in that order, as
# PowerPoint's ThemeColorScheme enumerates them (indices 1..12). Values
# are plain VBA/COM "RGB()" longs (0xBBGGRR ordering).

$p = $ppt.ActivePresentation

$cs = $p.SlideMaster.Theme.ThemeColorScheme

$officeThemeColors = @(
    0,            # dk1      000000
    16777215,     # lt1      FFFFFF
    6968388,      # dk2      44546A
    15132391,     # lt2      E7E6E6
    13998939,     # accent1  5B9BD5
    3243501,      # accent2  ED7D31
    10855845,     # accent3  A5A5A5
    49407,        # accent4  FFC000
    12874308,     # accent5  4472C4
    4697456,      # accent6  70AD47
    12673797,     # hlink    0563C1
    7491477       # folHlink 954F72
)

for ($i = 1; $i -le $cs.Count; $i++) {
    $cs.Item($i).RGB = $officeThemeColors[$i - 1]
}
